# Update scripts with new TPM values for Cntn1-Ptprz1 LR pair sheet
# (FAPs/MuSCs sending x FAPs/MuSCs target; "ECs" target cluster dropped)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cntn1"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.1199133333333333
$ws.Range("H2").Value = 0.35974
$ws.Range("I2").Value = 0.3099390012751145
$ws.Range("J2").Value = 0.3099390012751145
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05372733333333333
$ws.Range("N2").Value = 0.161182
$ws.Range("O2").Value = 0.1072370469527173
$ws.Range("P2").Value = 0.1072370469527173
$ws.Range("Q2").Value = 0.00644262363111111
$ws.Range("R2").Value = 0.05798361268
$ws.Range("S2").Value = 0.03323694323221776
$ws.Range("T2").Value = 0.03323694323221776

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cntn1"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1199133333333333
$ws.Range("H3").Value = 0.35974
$ws.Range("I3").Value = 0.3099390012751145
$ws.Range("J3").Value = 0.3099390012751145
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.4472873333333333
$ws.Range("N3").Value = 1.341862
$ws.Range("O3").Value = 0.8927629530472828
$ws.Range("P3").Value = 0.8927629530472827
$ws.Range("Q3").Value = 0.05363571509777777
$ws.Range("R3").Value = 0.4827214358799999
$ws.Range("S3").Value = 0.2767020580428968
$ws.Range("T3").Value = 0.2767020580428968

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cntn1"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2669800000000001
$ws.Range("H4").Value = 0.8009400000000001
$ws.Range("I4").Value = 0.6900609987248855
$ws.Range("J4").Value = 0.6900609987248854
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05372733333333333
$ws.Range("N4").Value = 0.161182
$ws.Range("O4").Value = 0.1072370469527173
$ws.Range("P4").Value = 0.1072370469527173
$ws.Range("Q4").Value = 0.01434412345333333
$ws.Range("R4").Value = 0.12909711108
$ws.Range("S4").Value = 0.07400010372049952
$ws.Range("T4").Value = 0.07400010372049952

# Row 5
$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Cntn1"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.2669800000000001
$ws.Range("H5").Value = 0.8009400000000001
$ws.Range("I5").Value = 0.6900609987248855
$ws.Range("J5").Value = 0.6900609987248854
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4472873333333333
$ws.Range("N5").Value = 1.341862
$ws.Range("O5").Value = 0.8927629530472828
$ws.Range("P5").Value = 0.8927629530472827
$ws.Range("Q5").Value = 0.1194167722533333
$ws.Range("R5").Value = 1.07475095028
$ws.Range("S5").Value = 0.616060895004386
$ws.Range("T5").Value = 0.6160608950043859

# Remove now-obsolete rows 6 and 7 (target cluster "ECs" combos dropped)
$ws.Rows("6:7").Delete()
